$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (values are read from source row's
# D/H/J/K/L/M/N/O/P columns and written into the target row's same columns).
# Columns A,B,C,E,F,G,I,Q,R are left untouched (they are constant across rows
# in this dataset, so the net effect is a per-row permutation of the
# Fecha/Variedad/Volumen/Precio.../Unidad/Origen/Precio $/Kg fields).
$rowMap = @{
    2 = 13
    3 = 14
    4 = 12
    5 = 15
    6 = 36
    7 = 29
    8 = 16
    9 = 6
    10 = 27
    11 = 10
    12 = 19
    13 = 8
    14 = 5
    15 = 31
    16 = 20
    17 = 32
    18 = 21
    19 = 18
    20 = 11
    21 = 23
    22 = 33
    23 = 35
    24 = 22
    25 = 2
    26 = 26
    27 = 30
    28 = 28
    29 = 24
    30 = 9
    31 = 7
    32 = 25
    33 = 4
    34 = 34
    35 = 17
    36 = 3
}

$cols = @("D","H","J","K","L","M","N","O","P")

# Snapshot every source cell's value BEFORE writing anything, since several
# rows both give to and receive from other rows (a true permutation) and a
# naive in-place copy would clobber a value before it is read.
$buffer = @{}
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range($col + $srcRow).Value()
    }
    $buffer[$destRow] = $rowVals
}

# Now write the buffered values back into their destination rows.
foreach ($destRow in $buffer.Keys) {
    $rowVals = $buffer[$destRow]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value = $rowVals[$col]
    }
}
